$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated word count for B2 ("Words"); formulas that depend on it
# (B4, D9, E9, D11, E11) recalculate automatically.
$ws.Range("B2").Value = 14169

# New blank cell styled like B2 but without the border (Times New Roman,
# no border) -- mirrors the new cellXfs entry added to the sheet.
$ws.Range("B9").Font.Name = "Times New Roman"

# New progress figure in B11, using the same look as B2 (Times New Roman,
# bordered) -- copy B2's formatting across rather than re-deriving the
# font/border so no stray style records get created.
$ws.Range("B11").Value = 14028
$ws.Range("B2").Copy()
$ws.Range("B11").PasteSpecial(-4122)

# Selection as last left by the user.
$ws.Range("G32").Select()
